$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 77 ("「諸都市の母、メッカ」" post).
# All rows below it (78..234) shift up by one automatically.
$ws.Rows.Item(77).Delete()
